$d = $word.ActiveDocument
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "Mode -") {
        $p.Range.Delete()
        break
    }
}
